$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
